# Update correlation_analysis results after fixing openjml handling.
# Touches two sheets:
#  - "all_tools"  (Worksheets.Item(1))  rows 9-12, plus a column-width swap (I<->J)
#  - "openjml"    (Worksheets.Item(5))  rows 9-12, plus a column-width shrink on J

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: all_tools
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("all_tools")

# Swap the widths stored for columns I (9) and J (10). We read the live
# values back from the model and cross-assign them (rather than writing new
# literals) so the swap stays internally consistent.
$wI = $ws.Columns.Item(9).ColumnWidth
$wJ = $ws.Columns.Item(10).ColumnWidth
$ws.Columns.Item(9).ColumnWidth = $wJ
$ws.Columns.Item(10).ColumnWidth = $wI

# Row 9
$ws.Range("G9").Value = 1143
$ws.Range("I9").Value = -0.166102596545867
$ws.Range("J9").Value = 0.01727604806480851
$ws.Range("K9").Value = -0.2527623213330977
$ws.Range("L9").Value = 0.01117475265921138

# Row 10
$ws.Range("G10").Value = 859
$ws.Range("I10").Value = -0.004414751593059719
$ws.Range("J10").Value = 0.9661430518463994
$ws.Range("K10").Value = 0.01233315619210278
$ws.Range("L10").Value = 0.932256616388448

# Row 11
$ws.Range("G11").Value = 859
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0.01007744778718968
$ws.Range("L11").Value = 0.9446252971798705

# Row 12
$ws.Range("G12").Value = 859
$ws.Range("I12").Value = 0.04719010357797921
$ws.Range("J12").Value = 0.632959082865391
$ws.Range("K12").Value = 0.07620564634846749
$ws.Range("L12").Value = 0.5988937623118764

# ---------------------------------------------------------------------------
# Sheet: openjml
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("openjml")

# Column J (10) narrows by one character unit.
$ws2.Columns.Item(10).ColumnWidth = $ws2.Columns.Item(10).ColumnWidth - 1

# Row 9
$ws2.Range("G9").Value = 751
$ws2.Range("H9").Value = 100
$ws2.Range("I9").Value = -0.1230941982211943
$ws2.Range("J9").Value = 0.107046316578524
$ws2.Range("K9").Value = -0.1671214681732231
$ws2.Range("L9").Value = 0.09653017580355105

# Row 10
$ws2.Range("G10").Value = 215
$ws2.Range("H10").Value = 50
$ws2.Range("I10").Value = -0.1410673005708742
$ws2.Range("J10").Value = 0.1918480862059539
$ws2.Range("K10").Value = -0.1594915550278049
$ws2.Range("L10").Value = 0.2685745446816231

# Row 11
$ws2.Range("G11").Value = 215
$ws2.Range("H11").Value = 50
$ws2.Range("I11").Value = -0.1104149035826812
$ws2.Range("J11").Value = 0.285180706372306
$ws2.Range("K11").Value = -0.1439811999024739
$ws2.Range("L11").Value = 0.318492693335327

# Row 12
$ws2.Range("G12").Value = 215
$ws2.Range("H12").Value = 50
$ws2.Range("I12").Value = 0.1113864967082373
$ws2.Range("J12").Value = 0.2778961820094916
$ws2.Range("K12").Value = 0.1553560043181458
$ws2.Range("L12").Value = 0.2813439520692285
